# analyze static SSD and erasure coding test
#
# 1) "Availability" sheet: MTR (D column) values for the redundant groups
#    drop from 72 hours to 24 hours (row 7 / "ssd" group is untouched).
# 2) "HW Architecture" sheet: the bandwidth labels in column D are bumped
#    up ("0.2G" -> "1.4G", "1.25G" -> "12.5G").
# 3) Selection/active-sheet bookkeeping: the user ends up back on the
#    "HW Architecture" tab with H29 selected, having last looked at E14
#    on the "Availability" tab.
# Everything on "Manual Calculation" / "etc" is derived via formulas, so
# it recalculates automatically once the source cells above change.

$wb = $excel.ActiveWorkbook

$wsAvail = $wb.Worksheets.Item("Availability")
$wsAvail.Range("D3").Value = 24
$wsAvail.Range("D4").Value = 24
$wsAvail.Range("D5").Value = 24
$wsAvail.Range("D6").Value = 24
$wsAvail.Range("D8").Value = 24
$wsAvail.Range("D9").Value = 24

$wsHw = $wb.Worksheets.Item("HW Architecture")
$wsHw.Range("D12:D75").Value = "1.4G"
$wsHw.Range("D2:D5").Value = "12.5G"
$wsHw.Range("D8:D11").Value = "12.5G"
$wsHw.Range("D76:D79").Value = "12.5G"
$wsHw.Range("D82:D85").Value = "12.5G"

# Leave a "last looked at" selection behind on Availability...
$wsAvail.Activate()
$wsAvail.Range("E14").Select()

# ...then return to HW Architecture as the active/visible tab.
$wsHw.Activate()
$wsHw.Range("H29").Select()

$excel.Calculate()
